# Update "想去人数" (F column) counts that were regenerated for this
# gh-pages data refresh. The same underlying row appears in multiple
# sheets ("展览" / "演出" feed into the aggregated "全部类型" sheet),
# so every occurrence of a given row's count needs to be bumped.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$map1 = @{
    2  = 884
    3  = 1009
    4  = 789
    5  = 868
    6  = 446
    7  = 686
    9  = 1285
    10 = 714
    14 = 37
    15 = 972
    17 = 405
    18 = 374
    20 = 586
    21 = 148
    22 = 634
    23 = 35
    24 = 1001
    25 = 14
}
foreach ($row in $map1.Keys) {
    $ws1.Range("F$row").Value = $map1[$row]
}

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 56

# Sheet "全部类型" (all types, aggregated view)
$ws4 = $wb.Worksheets.Item("全部类型")
$map4 = @{
    4  = 884
    5  = 1009
    6  = 789
    7  = 868
    8  = 446
    9  = 686
    11 = 1285
    12 = 714
    19 = 37
    20 = 972
    23 = 405
    24 = 374
    27 = 56
    28 = 586
    33 = 148
    34 = 634
    35 = 35
    36 = 1001
    37 = 14
}
foreach ($row in $map4.Keys) {
    $ws4.Range("F$row").Value = $map4[$row]
}
